# Actualización automática 2025-07-15 10:55:08
# Update sales figures on the "CUMPLIMIENTO MENSUAL" sheet:
#   - D2 (VENTA for OTROS) increases from 3456.76 to 3592.67
#   - E2 (POR CUMPLIR) is recalculated as -D2
#   - D4 (TOTAL VENTA) is recalculated as the sum of D2:D3
#   - E4 (TOTAL POR CUMPLIR) is recalculated as C4 - D4
#   - F4 (TOTAL CUMPLIMIENTO) is recalculated as D4 / C4

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# New venta value for the OTROS group (row 2)
$ws.Range("D2").Value = 3592.67

# POR CUMPLIR is the negative of VENTA for row 2
$ws.Range("E2").Value = -3592.67

# Recomputed TOTAL row (row 4) that rolls up rows 2 and 3
$ws.Range("D4").Value = 4536.18
$ws.Range("E4").Value = 9187.16
$ws.Range("F4").Value = 0.3305448965047867
